$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (shifts existing rows 10-23 down to 11-24)
$ws.Rows(10).Insert()

# Row 10 (new): Objetivos: + Portuguese objectives text
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Objetivo GeralPermitir aos estudantes que compreendam os mecanismos de obtenção da influencia de diversos fatores (variáveis independentes de um processo) sobre as variáveis resposta (dependentes), através da análise multivariada.Objetivos EspecíficosSaber planejar e executar um experimento fatorial completo e fracionadoSaber analisar os resultados propondo a condição de melhor ajuste que otimiza os valores da variável resposta na região experimental estudadaDominar, pelo menos, um software comercial sobre o assuntoSaber modelar um processo, com base em dados empíricos'
$ws.Range("C10").Value = 'Objetivo GeralPermitir aos estudantes que compreendam os mecanismos de obtenção da influencia de diversos fatores (variáveis independentes de um processo) sobre as variáveis resposta (dependentes), através da análise multivariada.Objetivos EspecíficosSaber planejar e executar um experimento fatorial completo e fracionadoSaber analisar os resultados propondo a condição de melhor ajuste que otimiza os valores da variável resposta na região experimental estudadaDominar, pelo menos, um software comercial sobre o assuntoSaber modelar um processo, com base em dados empíricos'
$ws.Rows(10).RowHeight = 60

# Row 11 (was old row 10): fix label + English objectives text (was wrongly duplicated)
$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'General objectiveTo allow students to understand the mechanisms of obtaining the influence of several factors (independent variables of a process) on the response variables (dependent), through the multivariate analysis.Specific objectivesKnow how to plan and execute a complete and fractional factorial experimentKnowing to analyze the results proposing the condition of better fit that optimizes the values of the response variable in the studied experimental regionManage at least one commercial software on the subjectKnow how to model a process, based on empirical data'
$ws.Range("C11").Value = 'General objectiveTo allow students to understand the mechanisms of obtaining the influence of several factors (independent variables of a process) on the response variables (dependent), through the multivariate analysis.Specific objectivesKnow how to plan and execute a complete and fractional factorial experimentKnowing to analyze the results proposing the condition of better fit that optimizes the values of the response variable in the studied experimental regionManage at least one commercial software on the subjectKnow how to model a process, based on empirical data'

# Row 12 (was old row 11): becomes "Docentes responsaveis:" label only, clear B/C, remove custom height
$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()
$ws.Rows(12).AutoFit()

# Row 13 (was old row 12): drop the A label, set B/C to the docente text
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '5840535 - Messias Borges Silva'
$ws.Range("C13").Value = '5840535 - Messias Borges Silva'
# B13 is a brand-new cell (row 12 previously had no B cell); pick up the
# correct column-B formatting (wrap text, vertical-top) from a sibling cell
# instead of the default that a freshly-created cell would otherwise get.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 14 (was old row 13): Programa resumido text fix
$ws.Range("B14").Value = 'Introdução Experimentação convencional Experimentos Fatoriais completos Experimentos Fatoriais fracionados Análise de variância Metodologia de superfície de resposta Método de Taguchi'
$ws.Range("C14").Value = 'Introdução Experimentação convencional Experimentos Fatoriais completos Experimentos Fatoriais fracionados Análise de variância Metodologia de superfície de resposta Método de Taguchi'

# Row 16 (was old row 15): Programa text fix
$ws.Range("B16").Value = 'Introdução Experimentação convencional Experimentos Fatoriais completos 2k , Experimentos Fatoriais fracionados 2k-p, Método de Plackett Burman,  Análise de variância Metodologia de superfície de resposta, Método de Taguchi .'
$ws.Range("C16").Value = 'Introdução Experimentação convencional Experimentos Fatoriais completos 2k , Experimentos Fatoriais fracionados 2k-p, Método de Plackett Burman,  Análise de variância Metodologia de superfície de resposta, Método de Taguchi .'

# Row 19 (was old row 18): Metodo text fix
$ws.Range("B19").Value = '2 provas escritas'
$ws.Range("C19").Value = '2 provas escritas'

# Row 20 (was old row 19): Criterio text fix
$ws.Range("B20").Value = 'Serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. MF = (0,40*P1 + 0,40*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'
$ws.Range("C20").Value = 'Serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. MF = (0,40*P1 + 0,40*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'

# Row 21 (was old row 20): Norma de recuperacao text fix
$ws.Range("B21").Value = 'Uma provas escrita com conteúdo de todo o semestre. NF = (MF + PR)/2, onde PR é uma prova de recuperação'
$ws.Range("C21").Value = 'Uma provas escrita com conteúdo de todo o semestre. NF = (MF + PR)/2, onde PR é uma prova de recuperação'

# Row 22 (was old row 21): Bibliografia text fix
$ws.Range("B22").Value = '1. MONTGOMERY, D.C., Design and Analysis of Experiments, Wiley, 19912. BOX, G.E.; HUNTER, W.G.; HUNTER, J.S., Statistic for Experimenters, John Wiley & Sons, New York, 1978. 3. TAGUCHI, G.; WU, YU-IN., Introduction to off-Line Quality Control. Central Japan Quality Control Association. Meieki Nakamura-Ku Magaya, Japan, 1979. 4. BRUNS, R.E., Como Fazer Experimentos, Editora UNICAMP, 2010. 5. COX, D.R., Planning of Experiments, Wiley 1976. 6. COX, G.M.; COCHRAN, W.G., Experimental Desing. Wiley 1976. 7. SILVA M.B. et al, Design of Experiments-Applications, Editora Intech, 2013'
$ws.Range("C22").Value = '1. MONTGOMERY, D.C., Design and Analysis of Experiments, Wiley, 19912. BOX, G.E.; HUNTER, W.G.; HUNTER, J.S., Statistic for Experimenters, John Wiley & Sons, New York, 1978. 3. TAGUCHI, G.; WU, YU-IN., Introduction to off-Line Quality Control. Central Japan Quality Control Association. Meieki Nakamura-Ku Magaya, Japan, 1979. 4. BRUNS, R.E., Como Fazer Experimentos, Editora UNICAMP, 2010. 5. COX, D.R., Planning of Experiments, Wiley 1976. 6. COX, G.M.; COCHRAN, W.G., Experimental Desing. Wiley 1976. 7. SILVA M.B. et al, Design of Experiments-Applications, Editora Intech, 2013'

# Row 24 (was old row 23): Requisitos text value, strip trailing newline
$ws.Range("B24").Value = 'LOB1049 -  Estatística Multivariada  (Requisito fraco)'
$ws.Range("C24").Value = 'LOB1049 -  Estatística Multivariada  (Requisito fraco)'

